$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (B1:E1) change
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 values (B2:E2) change
$ws.Range("B2").Value = 6.3118189179278641
$ws.Range("C2").Value = 7.583308857097931
$ws.Range("D2").Value = 12.184128701107621
$ws.Range("E2").Value = 11.548387016975301

# Row 3 values (B3:E3) change
$ws.Range("B3").Value = 5.2817447298443208
$ws.Range("C3").Value = 8.20560062016256
$ws.Range("D3").Value = 8.0830055985159763
$ws.Range("E3").Value = 12.343011764612253

# Update the selection to match the new range used in the edit
$ws.Range("B1:E3").Select()
